$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "team omega" row (row 2) down to row 3 first, so we don't
# disturb its (unstyled) formatting, then fill in the new "moose" row at row 2.
$ws.Cells.Item(3, 1).Value = "team omega"
$ws.Cells.Item(3, 2).Value = 75
$ws.Cells.Item(3, 3).Value = "2025-03-07T14:15:50.897662"

# Row 2: moose
$ws.Cells.Item(2, 1).Value = "moose"
$ws.Cells.Item(2, 2).Value = 76
$ws.Cells.Item(2, 3).Value = "2025-03-10T15:18:54.271858"

# Rows 4-7: new entries appended after "team omega" (row 3)
$ws.Cells.Item(4, 1).Value = "test "
$ws.Cells.Item(4, 2).Value = 75
$ws.Cells.Item(4, 3).Value = "2025-03-10T13:27:48.108708"

$ws.Cells.Item(5, 1).Value = "test"
$ws.Cells.Item(5, 2).Value = 66
$ws.Cells.Item(5, 3).Value = "2025-03-10T14:53:48.905717"

$ws.Cells.Item(6, 1).Value = "march 10"
$ws.Cells.Item(6, 2).Value = 66
$ws.Cells.Item(6, 3).Value = "2025-03-10T13:42:36.068089"

$ws.Cells.Item(7, 1).Value = "elon i"
$ws.Cells.Item(7, 2).Value = 63
$ws.Cells.Item(7, 3).Value = "2025-03-10T14:12:41.020226"
